$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$demoWs = $wb.Worksheets.Item(2)

# --- Row 2: SchoolName (A2) regenerated with a new random-looking school name ---
$ws.Cells.Item(2, 1).Value = "FPK12School27503"

# --- Row 3: E3 regenerated with a new numeric-looking id, stored as text ---
# A plain Value assignment of a numeric-looking string gets auto-converted to a
# number by Excel, losing the shared-string type. Author the text value in a far
# away scratch cell (forced to text via NumberFormat "@"), copy only the *value*
# over with PasteSpecial so the destination keeps its own style (s="2"), then
# remove the scratch cell again so it leaves no trace in the used range.
$scratch = $ws.Cells.Item(500, 500)
$scratch.NumberFormat = "@"
$scratch.Value = "39917"
$scratch.Copy()
$ws.Cells.Item(3, 5).PasteSpecial(-4163)
$excel.CutCopyMode = 0
$scratch.Delete()

# --- Row 4: D4 (fpk12teacher literal) cleared out, E4 reset back to the literal "0" ---
$ws.Cells.Item(4, 4).ClearContents()
# Sheet "DEMO" already stores the literal text "0" in E4; reuse it so the copied
# value keeps its text type, then paste only the value to preserve E4's own style.
$demoWs.Cells.Item(4, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Row 5: D5 (fpk12student literal) cleared out, E5 reset back to the literal "0" ---
$ws.Cells.Item(5, 4).ClearContents()
$demoWs.Cells.Item(5, 5).Copy()
$ws.Cells.Item(5, 5).PasteSpecial(-4163)
$excel.CutCopyMode = 0
